{"js": "// Word Document edit: \"San Junipero\" paragraph \u2014 change\n//   \"...explorar aventuras, permanentemente ou temporariamente. San Junipero...\"\n// to\n//   \"...explorar aventuras, de forma permanente ou temporariamente. San Junipero...\"\n//\n// Strategy: locate the paragraph by its distinctive leading text, then scope a\n// search to just that paragraph (the word \"permanentemente\" also appears in the\n// next paragraph, so a document-wide search would be ambiguous) and replace the\n// \" permanentemente\" occurrence (leading space included) with \" de forma permanente\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"San Junipero \u00e9 o quarto epis\u00f3dio\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'San Junipero \u00e9 o quarto epis\u00f3dio...' paragraph.\");\n}\n\nconst hits = target.search(\" permanentemente\", { matchCase: true });\nhits.load(\"items/text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find ' permanentemente' inside the target paragraph.\");\n}\n\nhits.items[0].insertText(\" de forma permanente\", \"Replace\");\nawait context.sync();\n", "ps1": "# Word Document edit: \"San Junipero\" paragraph - change\n#   \"...explorar aventuras, permanentemente ou temporariamente. San Junipero...\"\n# to\n#   \"...explorar aventuras, de forma permanente ou temporariamente. San Junipero...\"\n#\n# The word \"permanentemente\" (with its leading space) also occurs in the NEXT\n# paragraph (\"...ficar em San Junipero permanentemente, enquanto Yorkie quer...\")\n# so the Find/Replace is scoped to the Range of the specific paragraph that\n# starts with \"San Junipero e o quarto episodio\" to avoid touching that other\n# occurrence.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*San Junipero \u00e9 o quarto epis\u00f3dio*\") {\n        $target = $p.Range\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'San Junipero \u00e9 o quarto epis\u00f3dio...' paragraph.\"\n}\n\n$find = $target.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" permanentemente\"\n$find.Replacement.Text = \" de forma permanente\"\n\n# 0=wdFindStop, 1=wdFindContinue, 2=wdFindAsk ; 1=wdReplaceOne\n$result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n\nif (-not $result) {\n    throw \"Find/Replace of ' permanentemente' failed inside the target paragraph.\"\n}\n"}
